$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they keep their
# original textual precision (matching the source feed formatting)
# instead of being coerced into floating point numbers.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D14", "D15", "D16", "D18", "D19", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values
$ws.Range("D2").Value = "29.287.45"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.901.95"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "326.51"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "0.4650"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.3928"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.07888"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "0.9902"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").Value = "22.01"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").Value = "1.889.25"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "5.752"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "0.06968"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "88.35"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "0.000009989"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Value = "17.09"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "29.293.62"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "5.322"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").Value = "11.12"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "155.83"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "19.42"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "6.015"
$ws.Range("E27").Value = "  +2.58%  "
$ws.Range("D28").Value = "118.54"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "1.911"
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("D30").Value = "0.09380"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").Value = "0.9091"
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").Value = "5.284"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").Value = "3.217"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "1.186"
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("D36").Value = "0.05789"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").Value = "7.740"
$ws.Range("E39").Value = "  -3.52%  "
$ws.Range("D40").Value = "0.5716"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "0.1786"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").Value = "9.763"
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").Value = "11.99"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "0.5353"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "2.200"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").Value = "0.07041"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").Value = "1.855"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "2.571"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").Value = "1.060"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "71.30"
$ws.Range("E51").Value = "  -0.54%  "
